$wb = $excel.ActiveWorkbook

# Rename sheets: ObjectRepository -> WebObjectRepository, MobileRepository -> ObjectRepository
$wb.Worksheets.Item("MobileRepository").Name = "ObjectRepository_tmp"
$wb.Worksheets.Item("ObjectRepository").Name = "WebObjectRepository"
$wb.Worksheets.Item("ObjectRepository_tmp").Name = "ObjectRepository"

# Set selections on specific sheets
$wsMobileProps = $wb.Worksheets.Item("MobileProperties")
$wsMobileProps.Range("F20").Select()

$wsObjRepo = $wb.Worksheets.Item("ObjectRepository")
$wsObjRepo.Activate()
$wsObjRepo.Range("F19").Select()
